{"js": "// Replace the 25 \"NNN\u00d7N=\" equation prompts in the practice table with the\n// new set of equations, keeping everything else (formatting, layout,\n// blank answer rows, header date line) unchanged.\nconst replacements = [\n  [\"586\u00d76=\", \"519\u00d77=\"],\n  [\"257\u00d78=\", \"168\u00d78=\"],\n  [\"336\u00d73=\", \"256\u00d73=\"],\n  [\"808\u00d77=\", \"545\u00d75=\"],\n  [\"736\u00d79=\", \"760\u00d73=\"],\n  [\"428\u00d74=\", \"973\u00d74=\"],\n  [\"954\u00d79=\", \"631\u00d76=\"],\n  [\"115\u00d76=\", \"334\u00d77=\"],\n  [\"216\u00d72=\", \"286\u00d76=\"],\n  [\"302\u00d77=\", \"320\u00d79=\"],\n  [\"458\u00d75=\", \"108\u00d72=\"],\n  [\"589\u00d76=\", \"473\u00d79=\"],\n  [\"289\u00d78=\", \"581\u00d76=\"],\n  [\"833\u00d76=\", \"850\u00d78=\"],\n  [\"374\u00d79=\", \"450\u00d73=\"],\n  [\"265\u00d72=\", \"630\u00d74=\"],\n  [\"913\u00d73=\", \"138\u00d79=\"],\n  [\"421\u00d73=\", \"268\u00d75=\"],\n  [\"395\u00d79=\", \"684\u00d74=\"],\n  [\"771\u00d75=\", \"406\u00d79=\"],\n  [\"251\u00d79=\", \"432\u00d72=\"],\n  [\"558\u00d78=\", \"187\u00d73=\"],\n  [\"544\u00d73=\", \"142\u00d75=\"],\n  [\"166\u00d77=\", \"669\u00d76=\"],\n  [\"964\u00d74=\", \"365\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"NNN\u00d7N=\" equation prompts in the practice table with the\n# new set of equations, keeping everything else (formatting, layout,\n# blank answer rows, header date line) unchanged.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"586\u00d76=\", \"519\u00d77=\"),\n    @(\"257\u00d78=\", \"168\u00d78=\"),\n    @(\"336\u00d73=\", \"256\u00d73=\"),\n    @(\"808\u00d77=\", \"545\u00d75=\"),\n    @(\"736\u00d79=\", \"760\u00d73=\"),\n    @(\"428\u00d74=\", \"973\u00d74=\"),\n    @(\"954\u00d79=\", \"631\u00d76=\"),\n    @(\"115\u00d76=\", \"334\u00d77=\"),\n    @(\"216\u00d72=\", \"286\u00d76=\"),\n    @(\"302\u00d77=\", \"320\u00d79=\"),\n    @(\"458\u00d75=\", \"108\u00d72=\"),\n    @(\"589\u00d76=\", \"473\u00d79=\"),\n    @(\"289\u00d78=\", \"581\u00d76=\"),\n    @(\"833\u00d76=\", \"850\u00d78=\"),\n    @(\"374\u00d79=\", \"450\u00d73=\"),\n    @(\"265\u00d72=\", \"630\u00d74=\"),\n    @(\"913\u00d73=\", \"138\u00d79=\"),\n    @(\"421\u00d73=\", \"268\u00d75=\"),\n    @(\"395\u00d79=\", \"684\u00d74=\"),\n    @(\"771\u00d75=\", \"406\u00d79=\"),\n    @(\"251\u00d79=\", \"432\u00d72=\"),\n    @(\"558\u00d78=\", \"187\u00d73=\"),\n    @(\"544\u00d73=\", \"142\u00d75=\"),\n    @(\"166\u00d77=\", \"669\u00d76=\"),\n    @(\"964\u00d74=\", \"365\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
